$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new column before column A; this shifts the existing
# columns (old A->B, B->C, C->D, D->E) along with their widths/content.
$ws.Columns("A:A").Insert()

# ---- Column widths ----
# B..E already carry the correct widths from the shift. Give the new
# column A its own width (closest value Excel's column-width model can
# represent to the authored 27.7109375 is 27.6666... i.e. ColumnWidth 26.75).
$ws.Columns("A:A").ColumnWidth = 26.75

# ---- Copy header/body formatting into the newly inserted column A cells
# so they reuse the same cell styles as their row instead of creating new
# style records (matches s="1" header row / s="2" Arial row 3 styling). ----
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B3").Copy()
$ws.Range("A3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B3").Copy()
$ws.Range("C3").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = $false

# ---- Header row (row 1) ----
$ws.Range("A1").Value = "Dataset"
$ws.Range("B1").Value = "Description"
$ws.Range("C1").Value = "Condition1"
$ws.Range("D1").Value = "Condition2"
$ws.Range("E1").Value = "Location"

# ---- Row 2 (Tendons data set) ----
$ws.Range("A2").Value = "Tendons"
$ws.Range("B2").Value = "Tendons vs PSC "
$ws.Range("C2").Value = "tendon"
$ws.Range("D2").Value = "PSC"
$ws.Range("E2").Value = "tables/tendon vs PSC.xlsx"

# ---- Row 3 (Cerebral Organoids data set) ----
$ws.Range("A3").Value = "Cerebral Organoids day 40"
$ws.Range("B3").Value = "Cerebral Organoids day 40 compared to hPSC and this is a very long wall of text. I write a few more words to make it even longer"
$ws.Range("C3").Value = "Cerebral Organoids day 40"
$ws.Range("D3").Value = "hPSC"
$ws.Range("E3").Value = "tables/CO day 40 vs hPSC.xlsx"

# ---- Selection / active cell ----
[void]$ws.Range("D1").Select()

# ---- Window position ----
$excel.ActiveWindow.Left = 2790
